# 977: Add GS to extract process and GS tab to example files
$wb = $excel.ActiveWorkbook

# Add a new worksheet "GS" after the last existing sheet (CMS)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "GS"

# Header row matching the other extract tabs
$headerRange = $newSheet.Range("A1:I1")
$headerRange.Font.Color = 0

$newSheet.Range("A1").Value = "Contact_ID"
$newSheet.Range("B1").Value = "Contact_Date"
$newSheet.Range("C1").Value = "Contact_Type_Code"
$newSheet.Range("D1").Value = "Contact_Type_Desc"
$newSheet.Range("E1").Value = "OM_Name"
$newSheet.Range("F1").Value = "OM_Key"
$newSheet.Range("G1").Value = "OM_Grade"
$newSheet.Range("H1").Value = "OM_Team_Key"
$newSheet.Range("I1").Value = "OM_Provider_Code"

[void]$headerRange.Select()
